$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @("Task", "Estimate", "Actual")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 2).Value = $headers[$i]
}

$data = @(
    @(1, "Setup Project", 4, 9.43),
    @(2, "Setup Scene manager", 2, 1.2),
    @(3, "Player Controls", 30, 36),
    @(4, "Tree Collisions", 10, 7.54),
    @(5, "Soldier Logic", 20, 7.23),
    @(6, "Hospital Logic", 15, 8.52),
    @(7, "Game Over/Win Screen", 20, 10.24),
    @(8, "Assets/sprites", 30, 26.5)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 1).Value = $data[$r][0]
    $ws.Cells.Item($row, 2).Value = $data[$r][1]
    $ws.Cells.Item($row, 3).Value = $data[$r][2]
    $ws.Cells.Item($row, 4).Value = $data[$r][3]
}

$ws.Columns.Item(2).ColumnWidth = 22.140625 - 0.8333333333333334
$ws.Columns.Item(3).ColumnWidth = 8.7109375 - 0.8333333333333334
$ws.Columns.Item(4).ColumnWidth = 13.85546875 - 0.8333333333333334

$ws.Range("G10").Select() | Out-Null
